# Weekly price-data refresh: a new daily record is inserted at row 413,
# pushing all subsequent records down by one row (old row 413 -> 414, ...,
# old row 522 -> 523).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(413).Insert()

$ws.Cells.Item(413, 1).Value  = 10
$ws.Cells.Item(413, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(413, 3).Value  = "La Araucanía"
$ws.Cells.Item(413, 4).Value  = 44736
$ws.Cells.Item(413, 5).Value  = 9
$ws.Cells.Item(413, 6).Value  = 100112043
$ws.Cells.Item(413, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(413, 8).Value  = "Sin especificar"
$ws.Cells.Item(413, 9).Value  = "Primera"
$ws.Cells.Item(413, 10).Value = 200
$ws.Cells.Item(413, 11).Value = 19000
$ws.Cells.Item(413, 12).Value = 20000
$ws.Cells.Item(413, 13).Value = 19500
$ws.Cells.Item(413, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(413, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(413, 16).Value = 325
$ws.Cells.Item(413, 17).Value = 60
$ws.Cells.Item(413, 18).Value = "Hortaliza"
